$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells (includes two row swaps: 23<->24, 32<->33)
# Column D (Price) and E (Volume) are text cells in the source data; force
# text number-format before assignment so numeric-looking strings (e.g. '1.00',
# '593.46') are not silently coerced into real numbers by Excel's COM layer.

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '67.143.46'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  -0.54%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.458.99'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -1.50%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '593.46'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  -1.18%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '179.36'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +1.35%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.612'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  +3.54%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +0.13%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '3.458.95'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -1.34%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.139'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +5.03%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '6.95'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  -2.94%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.429'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -0.69%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.059.83'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -1.25%  '

# Row 14
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +1.98%  '

# Row 15
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -0.47%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '67.149.15'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  -0.46%  '

# Row 17
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -1.73%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '3.462.01'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -0.73%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.19'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -2.11%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '14.13'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  -3.62%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '392.70'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -0.62%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '7.90'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -1.75%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'LEO'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.77'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  +0.85%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'Dai'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '0.999'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +0.12%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.537'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  -0.48%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '71.58'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -2.69%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.0000121'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  -1.69%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.33'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +0.50%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.174'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -3.36%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  +0.51%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.11'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -0.80%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '2.05'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  -1.24%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Fetch.AI'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.39'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  -2.39%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '23.46'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  -1.39%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '7.30'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  -1.50%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.999'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -0.09%  '

# Row 37
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -4.40%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '160.47'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -1.77%  '

# Row 39
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -0.37%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.80'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  +9.21%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.86'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -4.72%  '

# Row 42
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -1.27%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '6.73'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -4.76%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '26.05'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -0.99%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.0718'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  -2.34%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.753.22'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  -2.20%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '26.10'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -5.04%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '41.32'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -2.97%  '

# Row 49
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -1.40%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '323.90'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -4.79%  '

# Row 51
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  -3.86%  '
